$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 107; this shifts the existing rows 107-184
# down to 108-185, preserving all of their data/formatting intact.
$ws.Rows.Item(107).EntireRow.Insert()

# Populate the newly inserted row 107 with the new weekly record.
$ws.Cells.Item(107, 1).Value = 7
$ws.Cells.Item(107, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(107, 3).Value = "Ñuble"
$ws.Cells.Item(107, 4).Value = 44447
$ws.Cells.Item(107, 5).Value = 16
$ws.Cells.Item(107, 6).Value = 100114001
$ws.Cells.Item(107, 7).Value = "Papa"
$ws.Cells.Item(107, 8).Value = "Rodeo"
$ws.Cells.Item(107, 9).Value = "1a (guarda)"
$ws.Cells.Item(107, 10).Value = 300
$ws.Cells.Item(107, 11).Value = 7000
$ws.Cells.Item(107, 12).Value = 7500
$ws.Cells.Item(107, 13).Value = 7250
$ws.Cells.Item(107, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(107, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(107, 16).Value = 290
$ws.Cells.Item(107, 17).Value = 25
$ws.Cells.Item(107, 18).Value = "Hortaliza"

# Preserve the date style (s="2") on the new D107 cell, matching the
# other date cells in column D.
$ws.Cells.Item(107, 4).NumberFormat = $ws.Cells.Item(108, 4).NumberFormat
